$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 809
$ws.Range("I28").Value = 459.5
$ws.Range("J28").Value = 1391.5
$ws.Range("K28").Value = 459.5
$ws.Range("L28").Value = 1391.5
$ws.Range("M28").Value = 25.5
$ws.Range("N28").Value = -2361.5

$ws.Range("H62").Value = 1992.6364
$ws.Range("I62").Value = 2008.7778
$ws.Range("J62").Value = 1920
$ws.Range("K62").Value = 2008.7778
$ws.Range("L62").Value = 1920
$ws.Range("M62").Value = -1384.7778
$ws.Range("N62").Value = -3168

$ws.Range("H65").Value = 1992.6364
$ws.Range("I65").Value = 2008.7778
$ws.Range("J65").Value = 1920
$ws.Range("K65").Value = 10043.889
$ws.Range("L65").Value = 9600
$ws.Range("M65").Value = -6923.889000000001
$ws.Range("N65").Value = -15840

$ws.Range("H106").Value = 2589.5715
$ws.Range("I106").Value = 2589.5715
$ws.Range("K106").Value = 2589.5715
$ws.Range("M106").Value = -1958.5715

$ws.Range("H107").Value = 279.875
$ws.Range("J107").Value = 313.57144
$ws.Range("L107").Value = 313.57144
$ws.Range("N107").Value = -4153.57144

$ws.Range("H138").Value = 4278.607
$ws.Range("I138").Value = 2137
$ws.Range("J138").Value = 5664.353
$ws.Range("K138").Value = 6411
$ws.Range("L138").Value = 16993.059
$ws.Range("M138").Value = -1271
$ws.Range("N138").Value = -27273.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1776.2195
$ws.Range("I74").Value = 945.2273
$ws.Range("J74").Value = 2738.4211
$ws.Range("K74").Value = 945.2273
$ws.Range("L74").Value = 2738.4211
$ws.Range("M74").Value = -71.22730000000001
$ws.Range("N74").Value = -4486.4211

$ws.Range("H77").Value = 1776.2195
$ws.Range("I77").Value = 945.2273
$ws.Range("J77").Value = 2738.4211
$ws.Range("K77").Value = 4726.136500000001
$ws.Range("L77").Value = 13692.1055
$ws.Range("M77").Value = -358.1365000000005
$ws.Range("N77").Value = -22428.1055

$ws.Range("H132").Value = 8186.8
$ws.Range("I132").Value = 8100.154
$ws.Range("K132").Value = 24300.462
$ws.Range("M132").Value = -21770.462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2492.7222
$ws.Range("I134").Value = 2507.9333
$ws.Range("J134").Value = 2416.6667
$ws.Range("K134").Value = 7523.7999
$ws.Range("L134").Value = 7250.000100000001
$ws.Range("M134").Value = -4988.7999
$ws.Range("N134").Value = -12320.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3238.5715
$ws.Range("I132").Value = 3677.7693
$ws.Range("J132").Value = 2524.875
$ws.Range("K132").Value = 11033.3079
$ws.Range("L132").Value = 7574.625
$ws.Range("M132").Value = -8503.3079
$ws.Range("N132").Value = -12634.625

$ws.Range("H134").Value = 1433.5186
$ws.Range("I134").Value = 1227.4762
$ws.Range("J134").Value = 2154.6667
$ws.Range("K134").Value = 3682.4286
$ws.Range("L134").Value = 6464.000100000001
$ws.Range("M134").Value = -1147.4286
$ws.Range("N134").Value = -11534.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5332
$ws.Range("I3").Value = 2990
$ws.Range("J3").Value = 10016
$ws.Range("K3").Value = 8970
$ws.Range("L3").Value = 30048
$ws.Range("M3").Value = -8858
$ws.Range("N3").Value = -30272

$ws.Range("H5").Value = 7194.839
$ws.Range("I5").Value = 659.3182
$ws.Range("J5").Value = 23170.555
$ws.Range("K5").Value = 1977.9546
$ws.Range("L5").Value = 69511.66500000001
$ws.Range("M5").Value = -1865.9546
$ws.Range("N5").Value = -69735.66500000001

$ws.Range("H22").Value = 17166.666
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 50000
$ws.Range("K22").Value = 2250
$ws.Range("L22").Value = 150000
$ws.Range("M22").Value = -2081
$ws.Range("N22").Value = -150338

$ws.Range("H27").Value = 17166.666
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 50000
$ws.Range("K27").Value = 2250
$ws.Range("L27").Value = 150000
$ws.Range("M27").Value = -2148
$ws.Range("N27").Value = -150204

$ws.Range("H37").Value = 523655.9
$ws.Range("J37").Value = 523655.9
$ws.Range("L37").Value = 1570967.7
$ws.Range("N37").Value = -1571191.7

$ws.Range("H68").Value = 1878.9452
$ws.Range("I68").Value = 1316.32
$ws.Range("J68").Value = 2171.9792
$ws.Range("K68").Value = 3948.96
$ws.Range("L68").Value = 6515.937600000001
$ws.Range("M68").Value = -3137.96
$ws.Range("N68").Value = -8137.937600000001

$ws.Range("H71").Value = 1878.9452
$ws.Range("I71").Value = 1316.32
$ws.Range("J71").Value = 2171.9792
$ws.Range("K71").Value = 11846.88
$ws.Range("L71").Value = 19547.8128
$ws.Range("M71").Value = -7790.879999999999
$ws.Range("N71").Value = -27659.8128

$ws.Range("H82").Value = 2620
$ws.Range("I82").Value = 1986.6666
$ws.Range("K82").Value = 5959.9998
$ws.Range("M82").Value = -5553.9998

$ws.Range("H85").Value = 2620
$ws.Range("I85").Value = 1986.6666
$ws.Range("K85").Value = 5959.9998
$ws.Range("M85").Value = -4555.9998

$ws.Range("H126").Value = 3071.6667
$ws.Range("I126").Value = 2686
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8058
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3118
$ws.Range("N126").Value = -24880

$ws.Range("H131").Value = 1196.043
$ws.Range("J131").Value = 1165.6405
$ws.Range("L131").Value = 3496.9215
$ws.Range("N131").Value = -13576.9215

$ws.Range("H132").Value = 2650.8
$ws.Range("J132").Value = 2527.1428
$ws.Range("L132").Value = 22744.2852
$ws.Range("N132").Value = -27804.2852

$ws.Range("H135").Value = 7194.839
$ws.Range("I135").Value = 659.3182
$ws.Range("J135").Value = 23170.555
$ws.Range("K135").Value = 5933.8638
$ws.Range("L135").Value = 208534.995
$ws.Range("M135").Value = -3398.8638
$ws.Range("N135").Value = -213604.995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 9285
$ws.Range("J109").Value = 9285
$ws.Range("L109").Value = 9285
$ws.Range("N109").Value = -11365

$ws.Range("H135").Value = 51780
$ws.Range("J135").Value = 51780
$ws.Range("L135").Value = 51780
$ws.Range("N135").Value = -61920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15450

$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1212.25
$ws.Range("I122").Value = 1242.5714
$ws.Range("K122").Value = 3727.7142
$ws.Range("M122").Value = -1277.7142

$ws.Range("H132").Value = 5538.6313
$ws.Range("I132").Value = 5607.9414
$ws.Range("K132").Value = 16823.8242
$ws.Range("M132").Value = -14293.8242
